# The commit rewords two sentences in the body of the e-mail:
#   1. "...harm the value of our company name." -> "...harm the value of the company name."
#      (with "was made" and "the" coming out as their own runs, and a gramStart/gramEnd
#       proofing mark around "was made", as Word's grammar checker would leave behind)
#   2. "Our company embraces the freedom to express oneself.  ... protect our company's
#       intellectual property..." ->
#      "At Webcor Builders, we embrace the freedom to express ourselves.  ... protect the
#       company's intellectual property..."
#      (with a spellStart/spellEnd proofing mark around "Webcor")
#
# Find.Execute's replace collapses a paragraph's runs into one, which loses the run/
# proofErr granularity the diff shows, so instead each whole paragraph is located and its
# content is replaced in one shot via Range.InsertXML with the exact run/proofErr layout.

$d = $word.ActiveDocument

function Set-ParagraphXml($needle, $xmlPayload) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            $start = $p.Range.Start
            $end = $p.Range.End
            $r = $d.Range($start, $end)
            $r.InsertXML($xmlPayload)
            return $true
        }
    }
    return $false
}

$payload1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">On November 1, the automated monitoring system will be expanded to include employees’ personal blogs.  This decision </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>was made</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> in order to prevent messages that could possibly harm the value of </w:t></w:r><w:r><w:t>the</w:t></w:r><w:r><w:t xml:space="preserve"> company name. </w:t></w:r><w:r><w:t xml:space="preserve"> If you wish to know more about this change, please visit this </w:t></w:r><w:r><w:t>link</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>http://www.</w:t></w:r><w:r><w:t>employee.</w:t></w:r><w:r><w:t>webcor.com/</w:t></w:r><w:r><w:t>policy.html#electronicmedia</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$payload2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">At </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Webcor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Builders, we embrace</w:t></w:r><w:r><w:t xml:space="preserve"> the freedom to express </w:t></w:r><w:r><w:t>ourselves</w:t></w:r><w:r><w:t xml:space="preserve">.  However, as an employee, we also hold a responsibility to protect </w:t></w:r><w:r><w:t>the</w:t></w:r><w:r><w:t xml:space="preserve"> company’s intellectual property and public image.  </w:t></w:r><w:r><w:t>We appreciate your cooperation in</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>keeping the company’s interest in mind</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

Set-ParagraphXml "On November 1, the automated monitoring system" $payload1 | Out-Null
Set-ParagraphXml "Our company embraces the freedom" $payload2 | Out-Null
